# Sinusoidal level 4 spreadsheet -> Tide height data
# "fixed manual guessing of frequency": replace the placeholder sinusoid
# sample data with real tide-height measurements, rename the header, and
# make the header bold + centered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (shared string) ---
$ws.Range("A1").Value = "Tide height (ft) x hours after jan 1 12 am"

# --- New data table (hours after 1/1 12am, tide height ft) ---
$data = @(
    @(5.85,   -0.05),
    @(11.82,   8.19),
    @(18.5,   -0.61),
    @(24.35,   7.39),
    @(30.58,  -0.17),
    @(36.55,   8.18),
    @(43.17,  -0.65),
    @(49.07,   7.55),
    @(55.37,  -0.22),
    @(61.33,   8.05),
    @(67.88,  -0.62),
    @(73.83,   7.68),
    @(80.2,   -0.19),
    @(86.17,   7.79),
    @(92.67,  -0.51),
    @(98.65,   7.75),
    @(104.15, -0.12),
    @(110.05,  7.44),
    @(116.52, -0.33),
    @(122.53,  7.76)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Give every data cell (including the newly added rows 8:21, which
# previously did not exist) the same plain, theme-coloured font that the
# original rows 2:7 used.
$dataRange = $ws.Range("A2:B21")
$dataRange.Font.Bold = $false
$dataRange.Font.ThemeColor = 1

# Make the header bold + centered (new font + alignment in styles.xml).
# A1:B1 is merged, so only the top-left cell (A1) needs to carry the
# style - styling B1 too would materialize an (otherwise absent) B1 cell.
$headerRange = $ws.Range("A1")
$headerRange.Font.Bold = $true
$headerRange.Font.ThemeColor = 1
$headerRange.HorizontalAlignment = -4108   # xlCenter
